$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c3Value = @'
list(au_id = c("https://openalex.org/A5039150887", "https://openalex.org/A5042987183", "https://openalex.org/A5013336788", "https://openalex.org/A5055005087"), au_display_name = c("Pooja S. Jagadish", "Marcus J. Hutchinson", "Scott Kirkpatrick", "Olivia Y. Hung"), au_orcid = c("https://orcid.org/0000-0002-6167-0820", NA, NA, "https://orcid.org/0000-0001-5668-2483"), author_position = c("first", "middle", "middle", "last"), au_affiliation_raw = c("", "Department of Medicine, University of Arizona-Tucson", 
"Division of Cardiology, Department of Medicine", "Division of Cardiology, Department of Medicine"), institution_id = c(NA, "https://openalex.org/I138006243", NA, NA), institution_display_name = c(NA, "University of Arizona", NA, NA), institution_ror = c(NA, "https://ror.org/03m2x1q45", NA, NA), institution_country_code = c(NA, "US", NA, NA), institution_type = c(NA, "education", NA, NA), institution_lineage = c(NA, "https://openalex.org/I138006243, https://openalex.org/I3132585960", NA, NA))
'@

$d3Value = @'
The patient is a postmenopausal female in her mid-50s with a reported history of Crohn's disease who presented on 3 separate occasions with gastrointestinal symptoms that progressed to chest pain with elevated troponin and wall motion abnormalities on transthoracic echocardiogram. High-sensitivity troponin peak, leukocyte count, inflammatory markers, and echocardiographic description for each episode are presented in the Figure. During the first episode, she additionally developed respiratory failure. The transthoracic echocardiogram demonstrated a left ventricular ejection fraction of 20% with apical hypokinesis and basal sparing, consistent with classic takotsubo cardiomyopathy. She underwent urgent coronary angiography showing no occlusive disease or plaque rupture. She was started on lisinopril for heart failure. At follow-up, echocardiogram demonstrated normal left ventricular ejection fraction and wall motion, and her heart failure medications were discontinued. In contrast, episode 2 (4 months after the first) was notable for posterior ST-elevation myocardial infarction with a lower troponin peak. Transthoracic echocardiogram revealed an ejection fraction of 49% with basal hypokinesis and apical sparing, consistent with reverse takotsubo. Emergent coronary angiography was again normal. She began lisinopril, but bradycardia prevented beta-blocker initiation. Computed tomography of the abdomen and pelvis revealed terminal ileitis. Upper endoscopy and colonoscopy showed non-bleeding duodenal ulcers and diverticulosis of the entire colon. The terminal ilium was poorly visualized. At follow-up 3 months later, echocardiogram demonstrated normal ejection fraction and wall movement. Her lisinopril was discontinued. Six months later, the patient underwent a colonoscopy. Post-procedure, she developed severe nausea and vomiting, then substernal chest pain. Electrocardiogram was concerning for anterior ST-elevation myocardial infarction. Transthoracic echocardiogram revealed a left ventricular ejection fraction of 38% and global hypokinesis to akinesis, consistent with the rarest form of takotsubo cardiomyopathy. Given her history, coronary angiogram was deferred. The patient declined cardiac magnetic resonance imaging. She was started on steroids and lisinopril, and her symptoms improved over the next 24 hours. Repeat echocardiogram showed total resolution of wall motion abnormalities with recovery of ejection fraction. She has had no recurrence of takotsubo cardiomyopathy for 6 months. She plans to aggressively manage her gastrointestinal symptoms but questions whether she needs to remain on her heart failure medications indefinitely. Discussions are ongoing. Takotsubo cardiomyopathy, also known as stress cardiomyopathy, occurs in 1%-2% of patients who present with acute coronary syndrome.1Ghadri J-R Wittstein IS Prasad A et al.International expert consensus document on takotsubo syndrome (Part I): clinical characteristics, diagnostic criteria, and pathophysiology.Eur Heart J. 2018; 39: 2032-2046Crossref PubMed Scopus (814) Google Scholar,2Templin C Ghadri JR Diekmann J et al.Clinical features and outcomes of takotsubo (stress) cardiomyopathy.N Engl J Med. 2015; 373: 929-938Crossref PubMed Scopus (1530) Google Scholar The International Takotsubo Diagnostic Criteria requires new electrocardiogram abnormalities or elevated cardiac biomarker levels, transient left ventricular dysfunction with segmental wall motion abnormalities, and a “trigger” (emotional, physical, neurologic). Patients with takotsubo cardiomyopathy are predominantly postmenopausal females and generally have negative evaluations for occlusive coronary atherosclerosis, pheochromocytoma, and myocarditis—although there are rare instances of these conditions co-existing.1Ghadri J-R Wittstein IS Prasad A et al.International expert consensus document on takotsubo syndrome (Part I): clinical characteristics, diagnostic criteria, and pathophysiology.Eur Heart J. 2018; 39: 2032-2046Crossref PubMed Scopus (814) Google Scholar,2Templin C Ghadri JR Diekmann J et al.Clinical features and outcomes of takotsubo (stress) cardiomyopathy.N Engl J Med. 2015; 373: 929-938Crossref PubMed Scopus (1530) Google Scholar Multiple variants of takotsubo cardiomyopathy have been identified based on the distribution of wall motion abnormalities seen on echocardiography.2Templin C Ghadri JR Diekmann J et al.Clinical features and outcomes of takotsubo (stress) cardiomyopathy.N Engl J Med. 2015; 373: 929-938Crossref PubMed Scopus (1530) Google Scholar,3Scantlebury DC Prasad A Diagnosis of takotsubo cardiomyopathy.Circ J. 2014; 78: 2129-2139Crossref PubMed Scopus (183) Google Scholar Frequency is based on the International Takotsubo Registry:1Ghadri J-R Wittstein IS Prasad A et al.International expert consensus document on takotsubo syndrome (Part I): clinical characteristics, diagnostic criteria, and pathophysiology.Eur Heart J. 2018; 39: 2032-2046Crossref PubMed Scopus (814) Google Scholar,3Scantlebury DC Prasad A Diagnosis of takotsubo cardiomyopathy.Circ J. 2014; 78: 2129-2139Crossref PubMed Scopus (183) Google Scholar,4Ghadri J-R Wittstein IS Prasad A et al.International expert consensus document on Takotsubo syndrome (Part II): diagnostic workup, outcome, and management.Eur Heart J. 2018; 39: 2047-2062Crossref PubMed Scopus (440) Google Scholar•Typical (81.7%) – Apical hypo- or dyskinesis with basal sparing•Mid-ventricular (14.6%) – Cuff-like area of hypokinesis in the mid-ventricle•“Reverse” (2.2%) – Basal hypokinesis with apical sparing•Focal (1.5%) – Areas of hypokinesis defined by the segment affected•Global (rare) – Hypokinesis in all segments The mechanism is not well understood but may be related to catecholamine surge, microvascular dysfunction, or coronary artery spasm.5Wittstein IS Thiemann DR Lima JAC et al.Neurohumoral features of myocardial stunning due to sudden emotional stress.N Engl J Med. 2005; 352: 539-548https://doi.org/10.1056/NEJMoa043046Crossref PubMed Scopus (2475) Google Scholar This is the first known case of recurrent takotsubo cardiomyopathy manifesting as 3 separate variants.
'@

$af3Value = @'
c("https://openalex.org/W2053449188", "https://openalex.org/W1913092753", "https://openalex.org/W4383710886", "https://openalex.org/W2889202876", "https://openalex.org/W2509841560", "https://openalex.org/W2463917843", "https://openalex.org/W2515818885", "https://openalex.org/W2035620702", "https://openalex.org/W3012475480", "https://openalex.org/W2464668688")
'@

$ws.Range("C3").Value = $c3Value
$ws.Range("D3").Value = $d3Value
$ws.Range("AF3").Value = $af3Value
